$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 25 and 26 swap their species-record values (columns A,B,E,F,G,H,L,Q,R)
# while the shared location/date/reporter columns remain unchanged.

$cols = @("A","B","E","F","G","H","Q","R")

foreach ($col in $cols) {
    $addr25 = "$col" + "25"
    $addr26 = "$col" + "26"
    $v25 = $ws.Range($addr25).Value2
    $v26 = $ws.Range($addr26).Value2
    $ws.Range($addr25).Value2 = $v26
    $ws.Range($addr26).Value2 = $v25
}

# L25 (empty cell) moves to L26; L25 ends up with no cell at all
$ws.Range("L25").Copy($ws.Range("L26"))
$ws.Range("L25").ClearContents()
